# Fix typo: ", углубил знания и о UML." -> ", углубил знания о UML."
# i.e. remove the stray "и " right before "о", keeping the (hidden)
# _GoBack bookmark anchored immediately before "о" instead of after it.

$d = $word.ActiveDocument

# Locate the exact text span "и о UML" so we are robust to any minor
# offset differences, then work out absolute character positions from it.
$anchor = $d.Content
$anchor.Find.ClearFormatting()
$found = $anchor.Find.Execute("и о UML", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target phrase to fix."
}

$iStart = $anchor.Start          # position of "и"
$oStart = $iStart + 2            # "и" + " " -> position of "о"

# Re-seat the (hidden) _GoBack bookmark *first*, while "о" is still its
# own untouched run, so it again sits right before "о" instead of right
# after it (it previously sat right after "о", before " UML"). Adding a
# bookmark in the middle of a run splits that run in two at the
# bookmark, which reproduces the original two-run + bookmark layout
# without disturbing the "о" run's own identity.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$gbRange = $d.Range($oStart, $oStart)
$d.Bookmarks.Add("_GoBack", $gbRange)

# Now remove the stray "и " (the letter and the following space) that
# precedes "о"; this merges the two (now adjacent, identically
# formatted) pieces of the preceding run back into one run, same as
# Word's own run coalescing on a text edit.
$toDelete = $d.Range($iStart, $oStart)
$toDelete.Delete()
